# Daily attendance processing - 2026-01-18 20:35:33
# Swap the order of names in the "Recorded By" column (G) so that entries
# reading "dnasr281@gmail.com, System" become "System, dnasr281@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
